# The deck has a row of SVM "Test Results w/ Different Parameters" slides,
# each with a bar chart whose value axis is pinned to a 50-100 scale. The
# chart on the "Digit Classification (LibSVM)" slide was missing its
# explicit axis Maximum (it only had Minimum=50), so its bars rendered on
# an inconsistent auto-scaled axis compared to its sibling charts. Fix it
# by explicitly setting the value axis Maximum back to 100, matching the
# other SVM charts.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(22)
$shp = $s.Shapes.Item(2)
$chart = $shp.Chart

# Axes(2,1) = the primary (xlPrimary) value axis (xlValue = 2)
$valueAxis = $chart.Axes(2, 1)
$valueAxis.MaximumScale = 100.0
